# Commit: "Added an FR for search"
#
# Adds a new Functional Requirement (FR) describing a tool-post search
# feature:
#   - "FRs" sheet        -> new row 47, numbered FR44 in that sheet's table
#   - "FRs Cleaning" sheet-> new row 46, numbered FR37 in the cleaned table
#     (rows 45-51, which only held leftover sequence numbers 37-43 in
#     column A, are removed and replaced by the single fully-populated
#     row 46)
#
# The editor finished their session with the "FRs Cleaning" sheet active.

$wb = $excel.ActiveWorkbook

$wsFRs = $wb.Worksheets.Item("FRs")
$wsCleaning = $wb.Worksheets.Item("FRs Cleaning")
$wsNFRs = $wb.Worksheets.Item("NFRs")

$searchDescription = "A search bar will be displayed.`nWhen the user searches, a list of tool-posts that match the search query will be displayed"
$searchRequirement = "The system must allow the user to search for a tool-post."
$createdNote = "Created 9/3/21"

# --- "FRs Cleaning" sheet: remove the leftover placeholder rows (45-51) ---
$wsCleaning.Range("A45:A51").EntireRow.Delete()

# --- Fill in the new row in "FRs Cleaning" (row 46) ---
# Order matches the authoring sequence: description & requirement text were
# typed here first, before the id cells on either sheet were filled in.
$wsCleaning.Range("D46").Value = $searchDescription
$wsCleaning.Range("C46").Value = $searchRequirement

# --- Fill in the new row in "FRs" (row 47) ---
$wsFRs.Range("D47").Value = $searchDescription
$wsFRs.Range("C47").Value = $searchRequirement
$wsFRs.Range("H47").Value = $createdNote
$wsFRs.Range("A47").Value = "FR44"

# --- Back to "FRs Cleaning" to finish the row ---
$wsCleaning.Range("H46").Value = $createdNote
$wsCleaning.Range("I46").Value = "U"
$wsCleaning.Range("A46").Value = "FR37"

$fCell = $wsCleaning.Range("F46")
$fCell.HorizontalAlignment = -4108
$fCell.VerticalAlignment = -4108
$fCell.WrapText = $true

$gCell = $wsCleaning.Range("G46")
$gCell.HorizontalAlignment = -4108
$gCell.VerticalAlignment = -4108
$gCell.WrapText = $true

# Row heights for the new wrapped-text rows.
$wsFRs.Rows.Item(47).RowHeight = 60
$wsCleaning.Rows.Item(46).RowHeight = 60

# --- Selections / active sheet, matching where the author ended up ---
$wsFRs.Activate()
$excel.ActiveWindow.ScrollRow = 37
$wsFRs.Range("A47:H47").Select()

$wsCleaning.Activate()
$wsCleaning.Range("K46").Select()

$wsNFRs.Select()
$wsCleaning.Activate()
